# "Main Assignment final commit"
# Sheet3 gets turned into a full user record (Name/Email/Password/Age),
# replacing the old dummy username/password placeholder row, the new
# email gets a mailto hyperlink, and Sheet3 becomes the active sheet.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")

# Header row: username/password -> Name/Email/Password/Age
$ws3.Range("A1").Value = "Name"
$ws3.Range("B1").Value = "Email"
$ws3.Range("C1").Value = "Password"
$ws3.Range("D1").Value = "Age"

# Data row: dummy/123456789 -> user1/user1@gmail.com/123456789/24
$ws3.Range("A2").Value = "user1"
$ws3.Range("B2").Value = "user1@gmail.com"
$ws3.Range("C2").Value = 123456789
$ws3.Range("D2").Value = 24

# Email cell becomes a mailto: hyperlink, styled like one
$ws3.Hyperlinks.Add($ws3.Range("B2"), "mailto:user1@gmail.com")
$ws3.Range("B2").Style = "Hyperlink"

# New Password column sized to fit its content
$ws3.Columns.Item(3).AutoFit()

# Sheet3 becomes the active/selected sheet (was Sheet2), with A7 selected
$ws3.Activate()
$ws3.Range("A7").Select()
